$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 120
$ws1.Range("F4").Value = 133
$ws1.Range("F5").Value = 1679
$ws1.Range("F6").Value = 1481
$ws1.Range("F7").Value = 272
$ws1.Range("F9").Value = 436

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 120
$ws4.Range("F4").Value = 133
$ws4.Range("F5").Value = 1679
$ws4.Range("F6").Value = 1481
$ws4.Range("F7").Value = 272
$ws4.Range("F10").Value = 436
